$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 44699
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = 100112001
$ws.Cells.Item(2, 7).Value = "Berenjena"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(2, 11).Value = 10000
$ws.Cells.Item(2, 12).Value = 10000
$ws.Cells.Item(2, 13).Value = 10000
$ws.Cells.Item(2, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(2, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(2, 16).Value = 167
$ws.Cells.Item(2, 17).Value = 60
$ws.Cells.Item(2, 18).Value = "Hortaliza"

$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
